$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.344.24"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "1.988.85"
$ws.Range("E3").Value = "  -4.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.017"
$ws.Range("E4").Value = "  +1.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.41"
$ws.Range("E5").Value = "  -4.19%  "
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4985"
$ws.Range("E7").Value = "  -4.73%  "
$ws.Range("E8").Value = "  -4.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.04"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08914"
$ws.Range("E10").Value = "  -4.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.110"
$ws.Range("E11").Value = "  -4.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.17"
$ws.Range("E12").Value = "  -6.48%  "
$ws.Range("D13").Value = "2.051.00"
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.950"
$ws.Range("E14").Value = "  -7.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.420"
$ws.Range("E15").Value = "  -7.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.017"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.33"
$ws.Range("E17").Value = "  -8.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001105"
$ws.Range("E18").Value = "  -4.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06705"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.43"
$ws.Range("E20").Value = "  -8.19%  "
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.912"
$ws.Range("E22").Value = "  -6.66%  "
$ws.Range("D23").Value = "29.435.39"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.94"
$ws.Range("E24").Value = "  -4.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.301"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.60"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.70"
$ws.Range("E27").Value = "  -5.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.234"
$ws.Range("E28").Value = "  -8.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.293"
$ws.Range("E29").Value = "  -8.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.43"
$ws.Range("E30").Value = "  -4.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.051"
$ws.Range("E31").Value = "  -7.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09890"
$ws.Range("E32").Value = "  -5.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.545"
$ws.Range("E33").Value = "  -6.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.811"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.804"
$ws.Range("E35").Value = "  -7.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02447"
$ws.Range("E36").Value = "  -7.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.201"
$ws.Range("E37").Value = "  -9.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.293"
$ws.Range("E38").Value = "  -3.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06341"
$ws.Range("E39").Value = "  -7.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6526"
$ws.Range("E40").Value = "  -6.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.55"
$ws.Range("E41").Value = "  -8.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2033"
$ws.Range("E42").Value = "  -8.08%  "
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6315"
$ws.Range("E44").Value = "  -7.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.42"
$ws.Range("E45").Value = "  -6.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.189"
$ws.Range("E46").Value = "  -6.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.301"
$ws.Range("E47").Value = "  -5.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.499"
$ws.Range("E48").Value = "  -3.89%  "
$ws.Range("E49").Value = "  -5.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06949"
$ws.Range("E50").Value = "  -4.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.125"
$ws.Range("E51").Value = "  -9.09%  "
